$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Paragraph 2: "dOur domain is data collected ... county+state names.d"
#   -> "Our domain is data collected ... county and state names."
# ---------------------------------------------------------------------------
$xmlPara2 = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t>Our domain is data collected where each data entry is associated with a county in the US. Data plugins simply read tabular data from many types of sources (</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>xls</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>, csv, web), and return a dictionary where sets of entries are indexed by county</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
  '<w:r><w:t>state</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> names.</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(2).Range.InsertXML($xmlPara2)

# ---------------------------------------------------------------------------
# Paragraph 3: "dOur framework then creates ... stdDev, etc.).d"
#   -> "Our framework then creates ... stdDev, etc.)."
# ---------------------------------------------------------------------------
$xmlPara3 = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t>Our framework then creates a new data-structure that associates each county with</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> geographic coordinates</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. This new data-structure also has library methods that return analysis of the given data in relation to time-frames, grouping of counties into bigger structures, and different types of analyses (sum, average, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>stdDev</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>, etc.).</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(3).Range.InsertXML($xmlPara3)

# ---------------------------------------------------------------------------
# Paragraph 4: "dThis new data-structure is then passed ... same entry.d"
#   -> "This new data-structure is then passed ... same entry."
# ---------------------------------------------------------------------------
$xmlPara4 = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t>This new data-structure is then passed to the display plug-in, which uses the data-structure' + [char]0x2019 + 's API to acquire summaries of the data in the form of raw columns</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">, where the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>ith</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> entry of a column is aligned with other columns to correspond to the same entry</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(4).Range.InsertXML($xmlPara4)

# ---------------------------------------------------------------------------
# Paragraph 5: "dPossible display plugins ... different times.d"
#   -> "Possible display plugins ... different times." (+ "sensitive" typo fix)
# ---------------------------------------------------------------------------
$xmlPara5 = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t>Possible</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> display plugins could be simple line plots, each line belonging to a county/state. </w:t></w:r>' +
  '<w:r><w:t>Similarly,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> for scatter plots, bubble charts, and even choropleths by using the coordinates</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> provided</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
  '<w:r><w:t>Given our time-</w:t></w:r>' +
  '<w:r><w:t>sensitive</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> data, display plug-ins could have a Google Earth functionality of displaying data with scroll-bar to change frames of data at different times.</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(5).Range.InsertXML($xmlPara5)

# Paragraph 6 (the bold "Describe your decisions..." bullet) is unchanged.

# ---------------------------------------------------------------------------
# Paragraph 7: "dWe decided to restrict ... would be easier. "
#   -> "We decided to restrict ... would be easier. " (drop leading d only)
# ---------------------------------------------------------------------------
$xmlPara7 = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t>We decided to restrict our specificity of the data to be data associated to a county in the US. We did this so we could create a hierarchy of membership amongst the inputted counties and their states, thus summarizing data per state</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> and even at the country level</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> would be easier.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(7).Range.InsertXML($xmlPara7)

# ---------------------------------------------------------------------------
# Bookmark "_GoBack" moves from the last paragraph (between the dotted-line
# run and "TO BE WRITTEN SOON") to the end of the "...use as necessary."
# paragraph. Delete the old bookmark first, then rebuild that paragraph
# (merging the "necessary." + stray "d" runs into one clean run) with the
# bookmark appended at the end.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$necessaryParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "The reusable functionality*") {
        $necessaryParaIndex = $i
        break
    }
}

$xmlNecessary = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:r><w:t xml:space="preserve">The reusable functionality the framework provides is the ability to organize very flexible data </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">in our domain </w:t></w:r>' +
  '<w:r><w:t>from any input supported by plugins and provide an object with a library having a multitude of data-analysis/filter tools to the display plugins to use as necessary.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$d.Paragraphs($necessaryParaIndex).Range.InsertXML($xmlNecessary)
